$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column cells hold plain-text values (e.g. "20.36", "0.5255")
# that look numeric, so force each edited cell to a text number format
# before writing the new value -- otherwise Excel would silently convert
# it into a floating point number instead of keeping the original text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.090.93'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.48'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5255'
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2611'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06352'
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.36'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07795'
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.506'
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.645.13'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5491'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅8237'
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.47'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.116.58'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.585'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '190.93'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '141.90'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1235'
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.238'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.432'
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05881'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.273'
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.531'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.263'
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9498'
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.785'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.410'
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5702'
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("E38").Value = '  +1.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.789'
$ws.Range("E39").Value = '  -2.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8474'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.84'
$ws.Range("E42").Value = '  +2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.025.70'
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.799.36'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9982'
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4299'
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.474'
$ws.Range("E48").Value = '  +2.59%  '
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.823'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09678'
$ws.Range("E51").Value = '  +0.02%  '
